$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two workers (row 16 / row 17) and update their "Salario Basico" values
# Row 16 becomes NORLIS LUCIA BANQUEZ LIDUEÑA (64698290), period 2301
$ws.Range("C16").Value = "64698290"
$ws.Range("D16").Value = "NORLIS LUCIA BANQUEZ LIDUEÑA"
$ws.Range("E16").Value = "2301"
$ws.Range("F16").Value = 39227
$ws.Range("G16").Value = 980657

# Row 17 becomes VICTOR MIRANDA TORRES (9297005), period 2302
$ws.Range("C17").Value = "9297005"
$ws.Range("D17").Value = "VICTOR MIRANDA TORRES"
$ws.Range("E17").Value = "2302"
$ws.Range("F17").Value = 46400
$ws.Range("G17").Value = 1423500
